$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="60.231.27"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Formula = '="  +3.63%  "'
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4163)
$ws.Range("D3").Formula = '="2.544.20"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Formula = '="  +3.60%  "'
$ws.Range("E3").Copy()
$ws.Range("E3").PasteSpecial(-4163)
$ws.Range("D4").Formula = '="1.00"'
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Formula = '="  +0.14%  "'
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial(-4163)
$ws.Range("D5").Formula = '="540.88"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Formula = '="  +2.52%  "'
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial(-4163)
$ws.Range("D6").Formula = '="145.28"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Formula = '="  +2.40%  "'
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial(-4163)
$ws.Range("E7").Formula = '="  -0.42%  "'
$ws.Range("E7").Copy()
$ws.Range("E7").PasteSpecial(-4163)
$ws.Range("E8").Formula = '="  +1.92%  "'
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial(-4163)
$ws.Range("D9").Formula = '="2.575.68"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Formula = '="  +4.44%  "'
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial(-4163)
$ws.Range("E10").Formula = '="  +3.81%  "'
$ws.Range("E10").Copy()
$ws.Range("E10").PasteSpecial(-4163)
$ws.Range("D11").Formula = '="0.160"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Formula = '="  +1.80%  "'
$ws.Range("E11").Copy()
$ws.Range("E11").PasteSpecial(-4163)
$ws.Range("D12").Formula = '="5.49"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Formula = '="  -0.36%  "'
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)
$ws.Range("D13").Formula = '="0.363"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Formula = '="  +5.61%  "'
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial(-4163)
$ws.Range("D14").Formula = '="2.995.01"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Formula = '="  +3.94%  "'
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial(-4163)
$ws.Range("D15").Formula = '="24.37"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Formula = '="  +3.79%  "'
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial(-4163)
$ws.Range("D16").Formula = '="60.186.07"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Formula = '="  +3.70%  "'
$ws.Range("E16").Copy()
$ws.Range("E16").PasteSpecial(-4163)
$ws.Range("D17").Formula = '="0.0000144"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Formula = '="  +6.72%  "'
$ws.Range("E17").Copy()
$ws.Range("E17").PasteSpecial(-4163)
$ws.Range("D18").Formula = '="2.579.56"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Formula = '="  +4.72%  "'
$ws.Range("E18").Copy()
$ws.Range("E18").PasteSpecial(-4163)
$ws.Range("D19").Formula = '="11.31"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Formula = '="  +1.70%  "'
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial(-4163)
$ws.Range("D20").Formula = '="4.36"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Formula = '="  +3.78%  "'
$ws.Range("E20").Copy()
$ws.Range("E20").PasteSpecial(-4163)
$ws.Range("D21").Formula = '="328.50"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Formula = '="  +3.58%  "'
$ws.Range("E21").Copy()
$ws.Range("E21").PasteSpecial(-4163)
$ws.Range("B22").Formula = '="Dai"'
$ws.Range("B22").Copy()
$ws.Range("B22").PasteSpecial(-4163)
$ws.Range("C22").Formula = '="https://coinranking.com/coin/MoTuySvg7+dai-dai"'
$ws.Range("C22").Copy()
$ws.Range("C22").PasteSpecial(-4163)
$ws.Range("D22").Formula = '="0.997"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Formula = '="  +0.01%  "'
$ws.Range("E22").Copy()
$ws.Range("E22").PasteSpecial(-4163)
$ws.Range("B23").Formula = '="Uniswap"'
$ws.Range("B23").Copy()
$ws.Range("B23").PasteSpecial(-4163)
$ws.Range("C23").Formula = '="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"'
$ws.Range("C23").Copy()
$ws.Range("C23").PasteSpecial(-4163)
$ws.Range("D23").Formula = '="5.96"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Formula = '="  +5.30%  "'
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial(-4163)
$ws.Range("D24").Formula = '="63.13"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Formula = '="  +5.41%  "'
$ws.Range("E24").Copy()
$ws.Range("E24").PasteSpecial(-4163)
$ws.Range("D25").Formula = '="0.437"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Formula = '="  +1.80%  "'
$ws.Range("E25").Copy()
$ws.Range("E25").PasteSpecial(-4163)
$ws.Range("E26").Formula = '="  +5.79%  "'
$ws.Range("E26").Copy()
$ws.Range("E26").PasteSpecial(-4163)
$ws.Range("D27").Formula = '="0.991"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Formula = '="  -1.25%  "'
$ws.Range("E27").Copy()
$ws.Range("E27").PasteSpecial(-4163)
$ws.Range("D28").Formula = '="8.04"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Formula = '="  +6.04%  "'
$ws.Range("E28").Copy()
$ws.Range("E28").PasteSpecial(-4163)
$ws.Range("D29").Formula = '="7.13"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Formula = '="  +5.59%  "'
$ws.Range("E29").Copy()
$ws.Range("E29").PasteSpecial(-4163)
$ws.Range("D30").Formula = '="0.0₃0806"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Formula = '="  +7.74%  "'
$ws.Range("E30").Copy()
$ws.Range("E30").PasteSpecial(-4163)
$ws.Range("E31").Formula = '="  +3.69%  "'
$ws.Range("E31").Copy()
$ws.Range("E31").PasteSpecial(-4163)
$ws.Range("D32").Formula = '="1.20"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Formula = '="  +0.37%  "'
$ws.Range("E32").Copy()
$ws.Range("E32").PasteSpecial(-4163)
$ws.Range("D33").Formula = '="164.80"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Formula = '="  +5.24%  "'
$ws.Range("E33").Copy()
$ws.Range("E33").PasteSpecial(-4163)
$ws.Range("E34").Formula = '="  +7.79%  "'
$ws.Range("E34").Copy()
$ws.Range("E34").PasteSpecial(-4163)
$ws.Range("E35").Formula = '="  -0.02%  "'
$ws.Range("E35").Copy()
$ws.Range("E35").PasteSpecial(-4163)
$ws.Range("D36").Formula = '="18.82"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Formula = '="  +3.00%  "'
$ws.Range("E36").Copy()
$ws.Range("E36").PasteSpecial(-4163)
$ws.Range("D37").Formula = '="4.47"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Formula = '="  +4.82%  "'
$ws.Range("E37").Copy()
$ws.Range("E37").PasteSpecial(-4163)
$ws.Range("D38").Formula = '="1.64"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Formula = '="  +4.55%  "'
$ws.Range("E38").Copy()
$ws.Range("E38").PasteSpecial(-4163)
$ws.Range("B39").Formula = '="OKB"'
$ws.Range("B39").Copy()
$ws.Range("B39").PasteSpecial(-4163)
$ws.Range("C39").Formula = '="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"'
$ws.Range("C39").Copy()
$ws.Range("C39").PasteSpecial(-4163)
$ws.Range("D39").Formula = '="37.19"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Formula = '="  +2.36%  "'
$ws.Range("E39").Copy()
$ws.Range("E39").PasteSpecial(-4163)
$ws.Range("B40").Formula = '="RenderToken"'
$ws.Range("B40").Copy()
$ws.Range("B40").PasteSpecial(-4163)
$ws.Range("C40").Formula = '="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"'
$ws.Range("C40").Copy()
$ws.Range("C40").PasteSpecial(-4163)
$ws.Range("D40").Formula = '="5.62"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Formula = '="  -0.96%  "'
$ws.Range("E40").Copy()
$ws.Range("E40").PasteSpecial(-4163)
$ws.Range("D41").Formula = '="302.58"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Formula = '="  +2.54%  "'
$ws.Range("E41").Copy()
$ws.Range("E41").PasteSpecial(-4163)
$ws.Range("B42").Formula = '="SuiNetwork"'
$ws.Range("B42").Copy()
$ws.Range("B42").PasteSpecial(-4163)
$ws.Range("C42").Formula = '="https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"'
$ws.Range("C42").Copy()
$ws.Range("C42").PasteSpecial(-4163)
$ws.Range("D42").Formula = '="0.840"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Formula = '="  +8.29%  "'
$ws.Range("E42").Copy()
$ws.Range("E42").PasteSpecial(-4163)
$ws.Range("B43").Formula = '="Filecoin"'
$ws.Range("B43").Copy()
$ws.Range("B43").PasteSpecial(-4163)
$ws.Range("C43").Formula = '="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"'
$ws.Range("C43").Copy()
$ws.Range("C43").PasteSpecial(-4163)
$ws.Range("D43").Formula = '="3.75"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Formula = '="  +4.40%  "'
$ws.Range("E43").Copy()
$ws.Range("E43").PasteSpecial(-4163)
$ws.Range("E44").Formula = '="  +4.91%  "'
$ws.Range("E44").Copy()
$ws.Range("E44").PasteSpecial(-4163)
$ws.Range("D45").Formula = '="0.993"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Formula = '="  -0.09%  "'
$ws.Range("E45").Copy()
$ws.Range("E45").PasteSpecial(-4163)
$ws.Range("D46").Formula = '="10.85"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Formula = '="  +0.89%  "'
$ws.Range("E46").Copy()
$ws.Range("E46").PasteSpecial(-4163)
$ws.Range("D47").Formula = '="126.98"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Formula = '="  +3.49%  "'
$ws.Range("E47").Copy()
$ws.Range("E47").PasteSpecial(-4163)
$ws.Range("D48").Formula = '="19.13"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Formula = '="  +5.45%  "'
$ws.Range("E48").Copy()
$ws.Range("E48").PasteSpecial(-4163)
$ws.Range("D49").Formula = '="0.0940"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Formula = '="  +2.80%  "'
$ws.Range("E49").Copy()
$ws.Range("E49").PasteSpecial(-4163)
$ws.Range("D50").Formula = '="0.0523"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Formula = '="  +2.98%  "'
$ws.Range("E50").Copy()
$ws.Range("E50").PasteSpecial(-4163)
$ws.Range("D51").Formula = '="0.0231"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Formula = '="  +3.29%  "'
$ws.Range("E51").Copy()
$ws.Range("E51").PasteSpecial(-4163)
